# Update "想去人数" (interested-count) values in column F for the
# "展览" and "全部类型" sheets, matching refreshed export data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — first data sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 145
$wsExpo.Range("F3").Value  = 1671
$wsExpo.Range("F4").Value  = 733
$wsExpo.Range("F6").Value  = 28
$wsExpo.Range("F7").Value  = 11763
$wsExpo.Range("F10").Value = 468
$wsExpo.Range("F11").Value = 391
$wsExpo.Range("F13").Value = 834
$wsExpo.Range("F14").Value = 13427
$wsExpo.Range("F15").Value = 13294
$wsExpo.Range("F20").Value = 264
$wsExpo.Range("F23").Value = 151

# Sheet "全部类型" (all types) — mirrors the same rows, but F7 lands on a
# slightly different refreshed value than the "展览" sheet.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 145
$wsAll.Range("F3").Value  = 1671
$wsAll.Range("F4").Value  = 733
$wsAll.Range("F6").Value  = 28
$wsAll.Range("F7").Value  = 11764
$wsAll.Range("F10").Value = 468
$wsAll.Range("F11").Value = 391
$wsAll.Range("F13").Value = 834
$wsAll.Range("F14").Value = 13427
$wsAll.Range("F15").Value = 13294
$wsAll.Range("F20").Value = 264
$wsAll.Range("F23").Value = 151
